$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the risk description text in H2 (this causes the shared-string
# table to drop the now-unused old string and renumber, matching the
# diff's cascading <v> index shifts across the rest of the sheet).
$ws.Range("H2").Value = "Understanding requirements, and lack of software fuctionalities "

# Move the active selection from G8 to H3.
$ws.Range("H3").Select()
